# Add "Rounded Notes" columns that split the previous single rounded-average
# column into two per-exam rounded columns plus a final rounded average
# column, on both worksheets.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Sheet 1 ("1ma1df01"): previously J held the rounded average of G,H.
# Now: J = round(G), K = round(H), L = round(average(J,K)).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Give the new header cells (J5, K6, L6) the same style as the existing
# bold/number-formatted header cell J6, then set their contents.
$ws1.Range("J6").Copy()
$ws1.Range("J5").PasteSpecial($xlPasteFormats)
$ws1.Range("K6").PasteSpecial($xlPasteFormats)
$ws1.Range("L6").PasteSpecial($xlPasteFormats)

$ws1.Range("J5").Value = "Rounded Notes"
$ws1.Range("J6").Value = "S1"
$ws1.Range("K6").Value = "S2"
$ws1.Range("L6").Value = "EOY"

for ($r = 7; $r -le 10; $r++) {
    $ws1.Range("J$r").Copy()
    $ws1.Range("K$r").PasteSpecial($xlPasteFormats)
    $ws1.Range("L$r").PasteSpecial($xlPasteFormats)

    $ws1.Range("J$r").Formula = "=ROUND(G$r*10)/10"
    $ws1.Range("K$r").Formula = "=ROUND(H$r*10)/10"
    $ws1.Range("L$r").Formula = "=ROUND(AVERAGE(J$r,K$r)*10)/10"
}

$ws1.Range("J7:J10").FormatConditions.Item(1).ModifyAppliesToRange($ws1.Range("J7:L10"))

# ---------------------------------------------------------------------------
# Sheet 2 ("2ma2dfb01"): previously I held the rounded average of F,G.
# Now: I = round(F), J = round(G), K = round(average(I,J)).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("I6").Copy()
$ws2.Range("I5").PasteSpecial($xlPasteFormats)
$ws2.Range("J6").PasteSpecial($xlPasteFormats)
$ws2.Range("K6").PasteSpecial($xlPasteFormats)

$ws2.Range("I5").Value = "Rounded Notes"
$ws2.Range("I6").Value = "S1"
$ws2.Range("J6").Value = "S2"
$ws2.Range("K6").Value = "EOY"

for ($r = 7; $r -le 29; $r++) {
    $ws2.Range("I$r").Copy()
    $ws2.Range("J$r").PasteSpecial($xlPasteFormats)
    $ws2.Range("K$r").PasteSpecial($xlPasteFormats)

    $ws2.Range("I$r").Formula = "=ROUND(F$r*10)/10"
    $ws2.Range("J$r").Formula = "=ROUND(G$r*10)/10"
    $ws2.Range("K$r").Formula = "=ROUND(AVERAGE(I$r,J$r)*10)/10"
}

$ws2.Range("I7:I29").FormatConditions.Item(1).ModifyAppliesToRange($ws2.Range("I7:K29"))
